# Tutorial 6 solution update:
#  - Date strings in column A switch from DD/MM/YYYY to DD-MM-YYYY.
#  - A handful of rows get their Real/Duplicate/Invalid/Absent flag
#    columns (D..H) recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the date separators in column A (rows 3-21) -------------------
# Some of the new strings (e.g. "01-08-2022") look like ambiguous
# day/month dates, so Excel would otherwise silently convert them into
# real date serial numbers. Prefixing with a leading apostrophe forces a
# text entry; resetting the Style back to "Normal" afterwards clears the
# "quote prefix" flag that operation leaves behind so the cell's styling
# stays identical to before.
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'" + $dates[$r]
    $cell.Style = "Normal"
}

# --- 2. Update the recomputed flag columns ---------------------------------
# Columns: D=Total Attendance Count, E=Real, F=Duplicate, G=Invalid, H=Absent
$ws.Cells.Item(3, 4).Value = 1    # D3
$ws.Cells.Item(3, 7).Value = 1    # G3

$ws.Cells.Item(4, 4).Value = 1    # D4
$ws.Cells.Item(4, 5).Value = 1    # E4
$ws.Cells.Item(4, 8).Value = 0    # H4

$ws.Cells.Item(5, 4).Value = 1    # D5
$ws.Cells.Item(5, 5).Value = 1    # E5
$ws.Cells.Item(5, 8).Value = 0    # H5

$ws.Cells.Item(12, 4).Value = 1   # D12
$ws.Cells.Item(12, 5).Value = 1   # E12
$ws.Cells.Item(12, 8).Value = 0   # H12

$ws.Cells.Item(13, 4).Value = 1   # D13
$ws.Cells.Item(13, 5).Value = 1   # E13
$ws.Cells.Item(13, 8).Value = 0   # H13
